$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the new "Mob_Skills" worksheet after the existing sheets
# ------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Mob_Skills"

# ------------------------------------------------------------------
# Data entry (column order chosen to reproduce the original authoring
# sequence: identifier/state/target/condition columns first, then the
# "Cancelable" column, then the header row, then the friendly names)
# ------------------------------------------------------------------

# Row 2 - Aquaring @ AL_HEAL (Heal)
$ws.Range("C2").Value = "Aquaring@AL_HEAL"
$ws.Range("D2").Value = "attack"
$ws.Range("K2").Value = "self"
$ws.Range("L2").Value = "myhpltmaxrate"

# Row 3 - Aquaring @ NPC_WATERATTACK (Water Attack)
$ws.Range("C3").Value = "Aquaring@NPC_WATERATTACK"
$ws.Range("D3").Value = "attack"
$ws.Range("K3").Value = "target"
$ws.Range("L3").Value = "always"

# Row 4 - Aquaring @ WZ_WATERBALL (Waterball)
$ws.Range("C4").Value = "Aquaring@WZ_WATERBALL"
$ws.Range("D4").Value = "attack"
$ws.Range("K4").Value = "target"
$ws.Range("L4").Value = "always"

# Cancelable column for all three rows
$ws.Range("J2").Value = "no"
$ws.Range("J3").Value = "no"
$ws.Range("J4").Value = "no"

# Header row
$ws.Range("B1").Value = "MobID"
$ws.Range("C1").Value = "Dummy value (info only)"
$ws.Range("D1").Value = "State"
$ws.Range("E1").Value = "SkillID"
$ws.Range("F1").Value = "SkillLv"
$ws.Range("G1").Value = "Rate"
$ws.Range("H1").Value = "CastTime"
$ws.Range("I1").Value = "Delay"
$ws.Range("J1").Value = "Cancelable"
$ws.Range("K1").Value = "Target"
$ws.Range("L1").Value = "Condition type"
$ws.Range("M1").Value = "Condition value"
$ws.Range("N1").Value = "val1"
$ws.Range("O1").Value = "val2"
$ws.Range("P1").Value = "val3"
$ws.Range("Q1").Value = "val4"
$ws.Range("R1").Value = "val5"
$ws.Range("S1").Value = "Emotion"
$ws.Range("T1").Value = "Chat"
$ws.Range("U1").Value = "Concat"

# Friendly skill names in column A (Waterball, then Heal, then Water Attack)
$ws.Range("A4").Value = "Waterball"
$ws.Range("A2").Value = "Heal"
$ws.Range("A3").Value = "Water Attack"

# ------------------------------------------------------------------
# Remaining numeric / formula cells for each row
# ------------------------------------------------------------------

# Row 2
$ws.Range("B2").Value = 3950
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 15
$ws.Range("G2").Formula = "=10000*0.1"
$ws.Range("H2").Formula = "=2.5*1000"
$ws.Range("I2").Formula = "=1000*2.5"
$ws.Range("M2").Value = 25
$ws.Range("U2").Formula = "=B2&"",""&C2&"",""&D2&"",""&E2&"",""&F2&"",""&G2&"",""&H2&"",""&I2&"",""&J2&"",""&K2&"",""&L2&"",""&M2&"",""&N2&"",""&O2&"",""&P2&"",""&Q2&"",""&R2&"",""&S2&"",""&T2"
$ws.Range("V2").Formula = "=""Casts level ""&F2&"" ""&A2&"" at ""&G2/100&""%""&"" on ""&PROPER(D2)&""."""

# Row 3
$ws.Range("B3").Value = 3950
$ws.Range("E3").Value = 184
$ws.Range("F3").Value = 5
$ws.Range("G3").Formula = "=10000*0.2"
$ws.Range("H3").Formula = "=0*1000"
$ws.Range("I3").Formula = "=1000*5"
$ws.Range("M3").Value = 0
$ws.Range("U3").Formula = "=B3&"",""&C3&"",""&D3&"",""&E3&"",""&F3&"",""&G3&"",""&H3&"",""&I3&"",""&J3&"",""&K3&"",""&L3&"",""&M3&"",""&N3&"",""&O3&"",""&P3&"",""&Q3&"",""&R3&"",""&S3&"",""&T3"
$ws.Range("V3").Formula = "=""Casts level ""&F3&"" ""&A3&"" at ""&G3/100&""%""&"" on ""&PROPER(D3)&""."""

# Row 4
$ws.Range("B4").Value = 3950
$ws.Range("E4").Value = 86
$ws.Range("F4").Value = 5
$ws.Range("G4").Formula = "=10000*0.03"
$ws.Range("H4").Formula = "=2*1000"
$ws.Range("I4").Formula = "=1000*0"
$ws.Range("M4").Value = 0
$ws.Range("U4").Formula = "=B4&"",""&C4&"",""&D4&"",""&E4&"",""&F4&"",""&G4&"",""&H4&"",""&I4&"",""&J4&"",""&K4&"",""&L4&"",""&M4&"",""&N4&"",""&O4&"",""&P4&"",""&Q4&"",""&R4&"",""&S4&"",""&T4"
$ws.Range("V4").Formula = "=""Casts level ""&F4&"" ""&A4&"" at ""&G4/100&""%""&"" on ""&PROPER(D4)&""."""

# ------------------------------------------------------------------
# Column width for column C (bestFit, matches the original author's
# manual widen of the skill-name column)
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 27.109375

# ------------------------------------------------------------------
# View state: Mob_Skills scrolled so column C is left-most, with
# V2:V4 selected (the computed description column)
# ------------------------------------------------------------------
$a1 = $ws.Activate()
$win = $wb.Application.ActiveWindow
$win.ScrollColumn = 3
$sel1 = $ws.Range("V2:V4").Select()

# ------------------------------------------------------------------
# Mob_Db sheet view: no longer the selected tab, scrolled right to
# column AJ, with BF2 selected
# ------------------------------------------------------------------
$wsDb = $wb.Worksheets.Item("Mob_Db")
$a2 = $wsDb.Activate()
$win2 = $wb.Application.ActiveWindow
$win2.ScrollColumn = 36
$sel2 = $wsDb.Range("BF2").Select()

# ------------------------------------------------------------------
# Mob_Skills is the active tab (rightmost sheet, index 2 / 0-based)
# ------------------------------------------------------------------
$a3 = $ws.Activate()
